# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Values for 展览 sheet (sheet1)
$ws1.Range("F2").Value = 157
$ws1.Range("F3").Value = 1780
$ws1.Range("F5").Value = 397
$ws1.Range("F11").Value = 22
$ws1.Range("F14").Value = 240
$ws1.Range("F16").Value = 34
$ws1.Range("F19").Value = 228
$ws1.Range("F21").Value = 439
$ws1.Range("F22").Value = 336
$ws1.Range("F27").Value = 733
$ws1.Range("F28").Value = 2532
$ws1.Range("F32").Value = 824
$ws1.Range("F36").Value = 376

# Values for 全部类型 sheet (sheet4) - duplicate data, note F19 differs (229 vs 228)
$ws4.Range("F2").Value = 157
$ws4.Range("F3").Value = 1780
$ws4.Range("F5").Value = 397
$ws4.Range("F11").Value = 22
$ws4.Range("F14").Value = 240
$ws4.Range("F16").Value = 34
$ws4.Range("F19").Value = 229
$ws4.Range("F21").Value = 439
$ws4.Range("F22").Value = 336
$ws4.Range("F27").Value = 733
$ws4.Range("F28").Value = 2532
$ws4.Range("F32").Value = 824
$ws4.Range("F36").Value = 376
